# LOQ4070.xlsx update
# - Fixes the "Objetivos:" row to show the real objectives text (was
#   pointing at the wrong shared string).
# - Inserts a new row for the "Docentes responsaveis" name (row 13),
#   pushing the remaining rows down by one.
# - Fills in previously-empty "Programa resumido", "Programa", "Metodo",
#   "Criterio", "Norma de recuperacao" and "Bibliografia" content cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Objetivos: row 10's B/C cells held the wrong text - put the real
#    course-objectives paragraph there instead.
# ---------------------------------------------------------------------
$objetivos = "Apresentar conceitos cinéticos e de fenômenos de transporte relativos a sistemas heterogêneos fluido-sólido com reações catalíticas bem como aplicações industriais de reatores heterogêneos catalíticos."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# ---------------------------------------------------------------------
# 2) Insert a row right before the old "Programa resumido:" row so the
#    responsible-professor name gets its own row (13), matching the
#    pattern used by "Objetivos:"/"Docentes responsaveis:" above it.
#    This shifts the former rows 13-23 down to 14-24.
# ---------------------------------------------------------------------
$ws.Rows("13:13").Insert()

$docente = "5963230 - Leandro Gonçalves de Aguiar"

$ws.Range("B13").Font.Bold = $false
$ws.Range("B13").WrapText = $true
$ws.Range("B13").VerticalAlignment = -4160
$ws.Range("B13").Value = $docente

$ws.Range("C13").Font.Bold = $false
$ws.Range("C13").Font.Color = 255
$ws.Range("C13").WrapText = $true
$ws.Range("C13").VerticalAlignment = -4160
$ws.Range("C13").Value = $docente

$ws.Range("A13").Clear()

# ---------------------------------------------------------------------
# 3) Fill in the syllabus / grading content that used to be blank.
# ---------------------------------------------------------------------
$programaResumido = @"
1. Conceitos gerais em catálise.
2. Tipos de sistemas catalíticos.
3. Interação fluido-sólido.
4. Velocidade das reações catalíticas gás-sólido.
5. Efeitos do transporte de massa e calor externo.
6. Transporte de massa interno.
7. Reatores heterogêneos catalíticos.
8. Modelos de reatores heterogêneos.
"@
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido
$ws.Rows("14:14").RowHeight = 60

$programa = @"
1. Conceitos gerais em catálise.
2. Tipos de sistemas catalíticos. Propriedades dos catalisadores sólidos.
3. Adsorção de um fluido sobre sólidos. Interação fluido-sólido.
4. Velocidade das reações catalíticas gás-sólido.
5. Efeitos do transporte de massa e calor externo.
6. Transporte de massa interno.
7. Reatores heterogêneos catalíticos.
8. Modelos de reatores heterogêneos.
"@
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

$metodo = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

$criterio = "Provas e trabalhos."
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

$norma = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("B21").Value = $norma
$ws.Range("C21").Value = $norma

$bibliografia = @"
FOGLER, H. S. Elementos de Engenharia das Reações Químicas. 3. ed. Rio de Janeiro: LTC Editora, 2002.
LEVENSPIEL, O. Chemical Reaction Engineering. 3. ed.  New York: John Wiley & Sons, 1998.
HILL, C.G. An Introduction to Chemical Engineering Kinetics and Reactor Design. New York: John Wiley&Sons, 1977.
SMITH, J.M. Chemical Engineering Kinetics. 3rd. ed. New York: McGraw-Hill, 1981.
DENBIGH, K.; TURNER, R. Introduction to Chemical Reaction Design. Cambridge: Cambridge University Press, 1970.
FROMENT, G.F.; Bischoff, K.B. Chemical Reactor Analysis and Design. 2nd. ed. New York: John Wiley & Sons, Inc. 1990.
Textos fornecidos pelo professor da disciplina
Artigos extraídos de revistas especializadas de Engenharia Química.
"@
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
